$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.188.66'
$ws.Range('E2').Value = '  +4.23%  '
$ws.Range('D3').Value = '3.324.99'
$ws.Range('E3').Value = '  +7.93%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '582.85'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').Value = '181.17'
$ws.Range('E6').Value = '  +6.69%  '
$ws.Range('D8').Value = '3.317.74'
$ws.Range('E8').Value = '  +7.82%  '
$ws.Range('D9').Value = '0.528'
$ws.Range('E9').Value = '  +3.58%  '
$ws.Range('D10').Value = '6.56'
$ws.Range('E10').Value = '  +2.59%  '
$ws.Range('D11').Value = '0.156'
$ws.Range('E11').Value = '  +3.55%  '
$ws.Range('D12').Value = '0.481'
$ws.Range('E12').Value = '  +2.33%  '
$ws.Range('D13').Value = '0.0000246'
$ws.Range('E13').Value = '  +2.87%  '
$ws.Range('D14').Value = '37.69'
$ws.Range('E14').Value = '  +5.24%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.792.35'
$ws.Range('E15').Value = '  +5.58%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '69.057.36'
$ws.Range('E16').Value = '  +4.16%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '0.123'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').Value = '7.29'
$ws.Range('E18').Value = '  +4.50%  '
$ws.Range('D19').Value = '3.190.46'
$ws.Range('E19').Value = '  +3.73%  '
$ws.Range('D20').Value = '16.86'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').Value = '498.04'
$ws.Range('E21').Value = '  +2.02%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '8.05'
$ws.Range('E22').Value = '  +4.00%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = '0.716'
$ws.Range('E23').Value = '  +4.09%  '
$ws.Range('D24').Value = '85.19'
$ws.Range('E24').Value = '  +2.97%  '
$ws.Range('D25').Value = '2.38'
$ws.Range('E25').Value = '  +7.75%  '
$ws.Range('D26').Value = '13.26'
$ws.Range('E26').Value = '  +4.40%  '
$ws.Range('D27').Value = '10.69'
$ws.Range('E27').Value = '  +5.37%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').Value = '8.19'
$ws.Range('E29').Value = '  +4.61%  '
$ws.Range('D30').Value = '2.43'
$ws.Range('E30').Value = '  +7.26%  '
$ws.Range('D31').Value = '2.67'
$ws.Range('E31').Value = '  +2.16%  '
$ws.Range('D32').Value = '29.48'
$ws.Range('E32').Value = '  +6.61%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').Value = '0.0₃0984'
$ws.Range('E33').Value = '  +8.11%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.115'
$ws.Range('E34').Value = '  +3.09%  '
$ws.Range('D35').Value = '0.990'
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('D36').Value = '5.93'
$ws.Range('E36').Value = '  +5.89%  '
$ws.Range('D37').Value = '0.991'
$ws.Range('E37').Value = '  +4.48%  '
$ws.Range('D38').Value = '48.33'
$ws.Range('E38').Value = '  +1.70%  '
$ws.Range('D39').Value = '0.329'
$ws.Range('E39').Value = '  +9.32%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '2.09'
$ws.Range('E40').Value = '  +5.83%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.128'
$ws.Range('E41').Value = '  +4.15%  '
$ws.Range('D42').Value = '49.90'
$ws.Range('E42').Value = '  +1.50%  '
$ws.Range('D43').Value = '8.52'
$ws.Range('E43').Value = '  +2.67%  '
$ws.Range('D44').Value = '2.78'
$ws.Range('E44').Value = '  +10.33%  '
$ws.Range('D45').Value = '409.70'
$ws.Range('E45').Value = '  +11.80%  '
$ws.Range('D46').Value = '2.878.14'
$ws.Range('E46').Value = '  +3.00%  '
$ws.Range('D47').Value = '27.46'
$ws.Range('E47').Value = '  +12.32%  '
$ws.Range('D48').Value = '0.0355'
$ws.Range('E48').Value = '  +2.36%  '
$ws.Range('D49').Value = '135.51'
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('D51').Value = '2.41'
$ws.Range('E51').Value = '  +11.42%  '
